$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.107.11'
$ws.Range("E2").Value = '  -2.48%  '
$ws.Range("D3").Value = '2.347.12'
$ws.Range("E3").Value = '  -3.34%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'310.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.13%  '
$ws.Range("D6").Value = "'85.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.16%  '
$ws.Range("E7").Value = '  -1.97%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -3.35%  '
$ws.Range("E10").Value = '  -3.01%  '
$ws.Range("D11").Value = "'30.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.90%  '
$ws.Range("D13").Value = '2.708.39'
$ws.Range("E13").Value = '  -3.38%  '
$ws.Range("D14").Value = "'6.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.50%  '
$ws.Range("D15").Value = "'14.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.66%  '
$ws.Range("D16").Value = '2.372.49'
$ws.Range("E16").Value = '  -2.08%  '
$ws.Range("D17").Value = "'0.761"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.60%  '
$ws.Range("D18").Value = '40.087.96'
$ws.Range("D19").Value = '0.0₃0903'
$ws.Range("E19").Value = '  -2.53%  '
$ws.Range("E20").Value = '  -2.79%  '
$ws.Range("D21").Value = "'68.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.82%  '
$ws.Range("D22").Value = "'10.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.77%  '
$ws.Range("D23").Value = "'235.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("E24").Value = '  -5.18%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  -3.18%  '
$ws.Range("D27").Value = "'23.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.02%  '
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("D29").Value = "'9.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.40%  '
$ws.Range("D30").Value = "'34.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").Value = "'153.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("E33").Value = '  -3.18%  '
$ws.Range("D34").Value = "'2.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.44%  '
$ws.Range("D35").Value = "'0.0719"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.11%  '
$ws.Range("E36").Value = '  -0.62%  '
$ws.Range("D37").Value = "'2.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.92%  '
$ws.Range("D38").Value = "'0.0991"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.27%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = "'1.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.94%  '
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").Value = "'15.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.70%  '
$ws.Range("D41").Value = "'3.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.11%  '
$ws.Range("D42").Value = '1.966.34'
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = "'0.0265"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.22%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'17.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.39%  '
$ws.Range("D46").Value = "'9.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.01%  '
$ws.Range("D47").Value = "'2.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.89%  '
$ws.Range("D48").Value = '2.567.96'
$ws.Range("D49").Value = "'93.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.07%  '
$ws.Range("D50").Value = "'70.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.50%  '
$ws.Range("D51").Value = "'50.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.68%  '
